$d = $word.ActiveDocument
$wNS = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function New-PkgXml($bodyInner) {
    return @"
<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='$wNS'><w:body>$bodyInner</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
}

# ---------------------------------------------------------------------------
# 1. Remove the "_GoBack" bookmark currently located after the "ColorBlind"
#    title run. It will be re-inserted later inside the new "Levels" section.
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------------
# 2. Genre section: delete the whole "Metroidvania" paragraph (merges away)
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Metroidvania`r") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 3. "Puzzle, shooter" paragraph -> split into two runs: "2d, p" / "uzzle, shooter"
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Puzzle, shooter`r") {
        $full = $d.Range($p.Range.Start, $p.Range.End - 1)
        $xml = New-PkgXml @"
<w:p>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>2d, p</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>uzzle, shooter</w:t></w:r>
</w:p>
"@
        $full.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------------
# 4. GREEN paragraph: "later screens" -> "next" + " screen"
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "GREEN:*") {
        $sr = $p.Range.Duplicate
        $sr.Find.Execute("later screens", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
        $tail = $d.Range($sr.Start, $p.Range.End - 1)
        $xml = New-PkgXml @"
<w:p>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>next</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> screen (moves left/right).</w:t></w:r>
</w:p>
"@
        $tail.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------------
# 5. RED paragraph: "all screens" -> "current and next screen" and keep the
#    trailing "colored death ground?" runs.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "RED:*") {
        $sr = $p.Range.Duplicate
        $sr.Find.Execute("all screens", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
        $tail = $d.Range($sr.Start, $p.Range.End - 1)
        $xml = New-PkgXml @"
<w:p>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>current and next screen</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>. Make </w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>colored death ground</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>?</w:t></w:r>
</w:p>
"@
        $tail.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------------
# 6. BLUE paragraph: "later screens" -> "the next" + " screens"
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "BLUE:*") {
        $sr = $p.Range.Duplicate
        $sr.Find.Execute("later screens", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
        $tail = $d.Range($sr.Start, $p.Range.End - 1)
        $xml = New-PkgXml @"
<w:p>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>the next</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> screens as well.</w:t></w:r>
</w:p>
"@
        $tail.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------------
# 7. Characters section "Enemies" paragraph: append " (red, blue, green,
#    black?)" and insert the new "Levels" section with 12 items plus the
#    "Zoom out big picture" closing line and the relocated "_GoBack" bookmark.
# ---------------------------------------------------------------------------
$enemiesParaIndex = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx++
    if ($p.Range.Text -eq "Enemies`r") {
        $enemiesParaIndex = $idx
    }
}
# the second "Enemies" paragraph (Characters section) is the one followed by "Roles"
$p = $d.Paragraphs($enemiesParaIndex)
$tail = $d.Range($p.Range.Start, $p.Range.End - 1)
$xml = New-PkgXml @"
<w:p>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Enemies</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> (red, blue, green, black?)</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:rPr><w:b/><w:u w:val='single'/><w:lang w:val='en-US'/></w:rPr></w:pPr>
<w:r><w:rPr><w:b/><w:u w:val='single'/><w:lang w:val='en-US'/></w:rPr><w:t>Levels</w:t></w:r>
</w:p>
<w:p>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>1. The story of life is an open canvas.</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> (Green paint for platform.)</w:t></w:r>
</w:p>
<w:p>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>2. Each choice makes changes to that canvas.</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> (Allow to choose between blue and red to deal with an enemy.)</w:t></w:r>
</w:p>
<w:p>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>3. Some changes are easier than others.</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> (Make it harder to get blue for the enemy. By needing to use a green to get it.)</w:t></w:r>
</w:p>
<w:p>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>4. Some might make our lives more difficult later on.</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> (Force to use red on enemy to kill it.)</w:t></w:r>
</w:p>
<w:p>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>5. Most of us are blind to each change as we make them.</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> (Just use some color platform puzzle.)</w:t></w:r>
</w:p>
<w:p>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:lastRenderedPageBreak/><w:t>6. Because changes most often are ones to the lens through which we see our canvas.</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> (Give a lot of different colors.)</w:t></w:r>
</w:p>
<w:p>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>7. And not the canvas truly there.</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> (Navigate over a pit with platforms</w:t></w:r>
<w:proofErr w:type='gramStart'/>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>. </w:t></w:r>
<w:proofErr w:type='gramEnd'/>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Stop the platforms with blue or kill enemies on them with red.)</w:t></w:r>
</w:p>
<w:p>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>8. Our life is a piece </w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>of that larger </w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>canvas &#8211; a piece of art.</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> (A lot of different colors.)</w:t></w:r>
</w:p>
<w:p>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>9. Even if small in comparison to the whole of the canvas.</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> (Create platforms to cross over a pit.)</w:t></w:r>
</w:p>
<w:p>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>10. We should treasure that small piece of art.</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> (Something </w:t></w:r>
<w:proofErr w:type='spellStart'/>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>something</w:t></w:r>
<w:proofErr w:type='spellEnd'/>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> enemies.)</w:t></w:r>
<w:bookmarkStart w:id='0' w:name='_GoBack'/>
<w:bookmarkEnd w:id='0'/>
</w:p>
<w:p>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>11. Whichever way we choose to paint it.</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> (</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Give lots of options to pass through stuff?</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>)</w:t></w:r>
</w:p>
<w:p>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>12. And strive to be happy with what we&#8217;ll have made at the end.</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> (Give colors </w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>to play around with </w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>and an end</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> jump to a hole of blackness.</w:t></w:r>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>)</w:t></w:r>
</w:p>
<w:p>
<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Zoom out big picture: As we all have just one canvas to paint.</w:t></w:r>
</w:p>
"@
$tail.InsertXML($xml)
